# Update the IP blacklist with new values and shrink the used range
# from A1:A23 to A1:A17 (rows 18-23 removed).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New IP values for rows 2-17 (row 9 keeps its original value 203.107.1.34)
$ws.Range("A2").Value2  = "112.126.94.107"
$ws.Range("A3").Value2  = "123.56.228.49"
$ws.Range("A4").Value2  = "139.99.236.139"
$ws.Range("A5").Value2  = "158.51.124.38"
$ws.Range("A6").Value2  = "158.51.124.56"
$ws.Range("A7").Value2  = "158.51.126.135"
$ws.Range("A8").Value2  = "185.215.113.66"
$ws.Range("A9").Value2  = "203.107.1.34"
$ws.Range("A10").Value2 = "217.8.117.10"
$ws.Range("A11").Value2 = "5.188.226.52"
$ws.Range("A12").Value2 = "66.187.4.127"
$ws.Range("A13").Value2 = "66.187.4.169"
$ws.Range("A14").Value2 = "66.187.4.92"
$ws.Range("A15").Value2 = "66.187.6.203"
$ws.Range("A16").Value2 = "92.63.197.112"
$ws.Range("A17").Value2 = "92.63.197.60"

# Remove the now-unused rows 18-23 so the used range shrinks back to A1:A17.
$ws.Range("A18:A23").ClearContents()
